$d = $word.ActiveDocument

# List of (old, new) text replacements, one per <w:t> run that changes.
# NOTE: Find.Execute(..., Replace=wdReplaceOne) always replaces the FIRST
# remaining occurrence of the search text anywhere in the document (it is
# not truly scoped to a sub-range in this runtime). Because a couple of the
# old values repeat and one of the new values ("51÷7=") happens to equal an
# old value that is still pending replacement later in the document, the
# operations below are ordered so that:
#   1) every old value is looked up (and consumed) in left-to-right/
#      top-to-bottom document order relative to *other* occurrences of the
#      same text, and
#   2) no newly-written text is created before any later pending search for
#      that same text has already run.
# This guarantees each Find.Execute call lands on the intended run.

$updates = @(
    @{old="2026-01-30 Friday"; new="2026-01-31 Saturday"},

    @{old="58÷5="; new="71÷2="},
    @{old="17÷3="; new="68÷8="},
    @{old="79÷4="; new="49÷7="},
    @{old="77÷3="; new="96÷4="},
    @{old="72÷3="; new="86÷2="},

    @{old="92÷2="; new="12÷3="},
    @{old="80÷8="; new="60÷9="},
    @{old="85÷3="; new="50÷9="},

    # Must run before the "16÷8=" -> "51÷7=" update below, otherwise that
    # update would create a "51÷7=" that this search would find instead of
    # the original one further down the document.
    @{old="51÷7="; new="82÷9="},

    @{old="16÷8="; new="51÷7="},
    @{old="84÷2="; new="59÷6="},

    @{old="15÷8="; new="80÷2="},
    @{old="40÷6="; new="32÷5="},
    @{old="66÷4="; new="71÷9="},
    @{old="54÷7="; new="65÷7="},
    @{old="38÷2="; new="37÷4="},

    @{old="57÷7="; new="68÷4="},
    @{old="91÷7="; new="30÷4="},
    @{old="21÷8="; new="43÷4="},
    @{old="17÷3="; new="80÷9="},
    @{old="66÷8="; new="21÷7="},

    @{old="87÷9="; new="65÷6="},
    @{old="84÷8="; new="78÷7="},
    @{old="87÷4="; new="17÷2="},
    @{old="14÷9="; new="15÷6="}
)

foreach ($u in $updates) {
    $rng = $d.Content
    $rng.Find.Execute($u.old, $true, $false, $false, $false, $false,
                       $true, 1, $false, $u.new, 1)
}
